# Highlight three bullet lines on the "Step 9" slide (Content Placeholder 2):
#   - "Implement the Read functions"                                              -> green  (00FF00)
#   - "Hint: Implement them in the order they appear in the data file"            -> yellow (FFFF00)
#   - "Using the debugger, step through each Read function, ..."                  -> red    (FF0000)

$p = $ppt.ActivePresentation

# Locate the slide that contains the target content placeholder (title "Step 9").
$targetSlide = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text.IndexOf("Implement the Read functions") -ge 0) {
                $targetSlide = $slide
            }
        }
    }
}

$shp = $null
for ($j = 1; $j -le $targetSlide.Shapes.Count; $j++) {
    $candidate = $targetSlide.Shapes.Item($j)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
        if ($candidate.TextFrame.TextRange.Text.IndexOf("Implement the Read functions") -ge 0) {
            $shp = $candidate
        }
    }
}

$tr = $shp.TextFrame.TextRange
$fullText = $tr.Text

function Set-RunHighlight($range, $fullText, $searchText, $rgb) {
    $idx = $fullText.IndexOf($searchText)
    if ($idx -ge 0) {
        $sub = $range.Characters($idx + 1, $searchText.Length)
        $sub.Font.Highlight.RGB = $rgb
    }
}

# RGB() packs as r + g*256 + b*65536 (matches VBA's RGB function / OOXML srgbClr).
Set-RunHighlight $tr $fullText "Implement the Read functions" 65280     # 00FF00
Set-RunHighlight $tr $fullText "Hint: Implement them in the order they appear in the data file" 65535   # FFFF00
Set-RunHighlight $tr $fullText "Using the debugger, step through each Read function, making sure that each object is created and assigned data correctly" 255   # FF0000
